$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D38").Value = 933.12
$ws1.Range("E38").Value = 193.91
$ws1.Range("M38").Value = 1229.62
$ws1.Range("M46").Value = 1606.1

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F38").Value = 3757.42
$ws2.Range("F46").Value = 1606.1
$ws2.Range("F57").Value = 23289.89

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D3").Value = 15535.01
$ws3.Range("E3").Value = -651.130000000001
$ws3.Range("F3").Value = 1.043747329325418

$ws3.Range("D4").Value = 193.91
$ws3.Range("E4").Value = 452.09
$ws3.Range("F4").Value = 0.3001702786377709

$ws3.Range("D15").Value = 4993.04
$ws3.Range("E15").Value = 15696.96
$ws3.Range("F15").Value = 0.2413262445625906

$ws3.Range("D16").Value = 14747.8
$ws3.Range("E16").Value = 43973.43000000001
$ws3.Range("F16").Value = 0.2511493713602388

$ws3.Range("D19").Value = 37045.54
$ws3.Range("E19").Value = 72823.21000000001
$ws3.Range("F19").Value = 0.3371799533534331

# Column E width change: 24 -> 23
# (The COM layer stores width as input+5/6, so feed it back-adjusted
#  so the serialized OOXML "width" attribute comes out to exactly 23.)
$ws3.Columns.Item(5).ColumnWidth = 22.166666666666668
